$wb = $excel.ActiveWorkbook

# Sheet "OFF" - Week 13 row (R, row 3) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 204
$wsOff.Range("C3").Value = 146
$wsOff.Range("D3").Value = 36
$wsOff.Range("E3").Value = 19
$wsOff.Range("G3").Value = 4

# Sheet "DEF" - Week 13 row (R, row 3) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 139
$wsDef.Range("C3").Value = 100
$wsDef.Range("D3").Value = 37
$wsDef.Range("E3").Value = 18
$wsDef.Range("G3").Value = 4
